# Change the household_id field in the household members "survey" form
# from a plain editable "text" field (with readonly=TRUE and a comments
# column) into a non-editable "note" field that just echoes the household
# id back to the user.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# The "readonly" and "comments" columns (F and G) are no longer needed once
# household_id becomes a note - Excel shifts the old column H (hideInContents)
# left into F automatically.
$ws.Columns("F:G").Delete()

# Row 2 is the household_id question. Turn it into a note:
#  - type: text -> note
#  - name: household_id -> (cleared, notes don't bind to a field name here)
#  - display.text: swap the old barcode-lookup hint for a note message
$ws.Range("A2").Value = "note"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "Data for household: {{household_id}}"

# Restore the cursor/selection to where it ends up after the edit.
$ws.Range("F6").Select()
